$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    # The source data is stored as text (inline strings) even though it
    # looks numeric, so force text formatting while writing so the value
    # is kept as literal text instead of being reinterpreted as a number,
    # then restore the default ("Normal") style so no stray number format
    # is left behind on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 (iteration 1): xi unchanged, f(xi) and Error updated
Set-TextValue "C2" "8.125"
Set-TextValue "D2" "0.552483908235682"

# Row 3 (iteration 2): xi, f(xi) and Error updated
Set-TextValue "B3" "0.947516091764318"
Set-TextValue "C3" "-1.60298414312629"
Set-TextValue "D3" "0.117187901308275"

# Row 4 (iteration 3): new row
Set-TextValue "A4" "3"
Set-TextValue "B4" "1.06470399307259"
Set-TextValue "C4" "-0.180331565072615"
Set-TextValue "D4" "0.0132080976575131"

# Row 5 (iteration 4): new row
Set-TextValue "A5" "4"
Set-TextValue "B5" "1.07791209073011"
Set-TextValue "C5" "-0.0019872259809312"
Set-TextValue "D5" "0.0001455571464379"

# Row 6 (iteration 5): new row
Set-TextValue "A6" "5"
Set-TextValue "B6" "1.07805764787654"
Set-TextValue "C6" "-2.36726886382339e-07"
Set-TextValue "D6" "1.73393919222775e-08"
